$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Names.Item("_xlchart.v1.2").Delete()
$wb.Names.Item("_xlchart.v1.3").Delete()

$ws.Range("D18").Value = "Mean increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Value = "Median increase"
$ws.Range("F18").Font.Bold = $true

$ws.Range("F19").Formula = "= ((E10 / 113.658804) * 100) - 100"
$ws.Range("D19").Formula = "= ((E3 / 114.202998) * 100) - 100"
$ws.Range("F19").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Range("D23").Select()
